$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 5020
$ws.Range("E2").Value = 141
$ws.Range("F2").Value = 141
$ws.Range("G2").Value = 46
$ws.Range("H2").Value = 38
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 37
$ws.Range("K2").Value = 5635
$ws.Range("L2").Value = 3949
$ws.Range("M2").Value = 1686
$ws.Range("N2").Value = 1395
$ws.Range("O2").Value = 291
$ws.Range("P2").Value = 108
$ws.Range("Q2").Value = 205
$ws.Range("R2").Value = -343
$ws.Range("S2").Value = 117
$ws.Range("T2").Value = 379
$ws.Range("U2").Value = -175
$ws.Range("V2").Value = 2935
$ws.Range("W2").Value = 2.81
$ws.Range("X2").Value = 0.76
$ws.Range("Y2").Value = 0.09
$ws.Range("Z2").Value = 0.68
$ws.Range("AA2").Value = 234.19
$ws.Range("AB2").Value = 1101.03
$ws.Range("AC2").Value = 5
$ws.Range("AD2").Value = 928.96
$ws.Range("AE2").Value = 5572
$ws.Range("AF2").Value = 0.82
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 25035871

# Row 3
$ws.Range("D3").Value = 5194
$ws.Range("E3").Value = 311
$ws.Range("F3").Value = 311
$ws.Range("G3").Value = 189
$ws.Range("H3").Value = 135
$ws.Range("I3").Value = 104
$ws.Range("J3").Value = 31
$ws.Range("K3").Value = 5858
$ws.Range("L3").Value = 4102
$ws.Range("M3").Value = 1756
$ws.Range("N3").Value = 1477
$ws.Range("O3").Value = 278
$ws.Range("P3").Value = 108
$ws.Range("Q3").Value = 382
$ws.Range("R3").Value = -296
$ws.Range("S3").Value = -59
$ws.Range("T3").Value = 324
$ws.Range("U3").Value = 58
$ws.Range("V3").Value = 2990
$ws.Range("W3").Value = 5.99
$ws.Range("X3").Value = 2.6
$ws.Range("Y3").Value = 7.22
$ws.Range("Z3").Value = 2.35
$ws.Range("AA3").Value = 233.69
$ws.Range("AB3").Value = 1201.56
$ws.Range("AC3").Value = 414
$ws.Range("AD3").Value = 17.51
$ws.Range("AE3").Value = 5952
$ws.Range("AF3").Value = 1.22
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 25035871

# Row 4
$ws.Range("D4").Value = 5575
$ws.Range("E4").Value = 348
$ws.Range("F4").Value = 348
$ws.Range("G4").Value = 96
$ws.Range("H4").Value = 84
$ws.Range("I4").Value = 44
$ws.Range("J4").Value = 39
$ws.Range("K4").Value = 7315
$ws.Range("L4").Value = 5332
$ws.Range("M4").Value = 1983
$ws.Range("N4").Value = 1666
$ws.Range("O4").Value = 317
$ws.Range("P4").Value = 108
$ws.Range("Q4").Value = 213
$ws.Range("R4").Value = -1180
$ws.Range("S4").Value = 1061
$ws.Range("T4").Value = 1118
$ws.Range("U4").Value = -906
$ws.Range("V4").Value = 3755
$ws.Range("W4").Value = 6.25
$ws.Range("X4").Value = 1.5
$ws.Range("Y4").Value = 2.83
$ws.Range("Z4").Value = 1.27
$ws.Range("AA4").Value = 268.84
$ws.Range("AB4").Value = 1380
$ws.Range("AC4").Value = 178
$ws.Range("AD4").Value = 47.2
$ws.Range("AE4").Value = 6711
$ws.Range("AF4").Value = 1.25
$ws.Range("AG4").Value = 43
$ws.Range("AH4").Value = 0.52
$ws.Range("AI4").Value = 24.17
$ws.Range("AJ4").Value = 25035871

# Row 5
$ws.Range("D5").Value = 6295
$ws.Range("E5").Value = 322
$ws.Range("F5").Value = 322
$ws.Range("G5").Value = 35
$ws.Range("H5").Value = 9
$ws.Range("I5").Value = -22
$ws.Range("J5").Value = 31
$ws.Range("K5").Value = 8275
$ws.Range("L5").Value = 6314
$ws.Range("M5").Value = 1961
$ws.Range("N5").Value = 1624
$ws.Range("O5").Value = 337
$ws.Range("P5").Value = 112
$ws.Range("Q5").Value = 278
$ws.Range("R5").Value = -1313
$ws.Range("S5").Value = 972
$ws.Range("T5").Value = 862
$ws.Range("U5").Value = -584
$ws.Range("V5").Value = 4771
$ws.Range("W5").Value = 5.12
$ws.Range("X5").Value = 0.15
$ws.Range("Y5").Value = -1.32
$ws.Range("Z5").Value = 0.12
$ws.Range("AA5").Value = 322.02
$ws.Range("AB5").Value = 1319.36
$ws.Range("AC5").Value = -85
$ws.Range("AD5").Value = -81.66
$ws.Range("AE5").Value = 6316
$ws.Range("AF5").Value = 1.1
$ws.Range("AG5").Value = 43
$ws.Range("AH5").Value = 0.62
$ws.Range("AI5").Value = -51.17
$ws.Range("AJ5").Value = 25921252

# Row 6
$ws.Range("D6").Value = 5882
$ws.Range("E6").Value = 226
$ws.Range("F6").Value = 226
$ws.Range("G6").Value = -32
$ws.Range("H6").Value = -70
$ws.Range("I6").Value = -77
$ws.Range("K6").Value = 8103
$ws.Range("L6").Value = 6238
$ws.Range("M6").Value = 1866
$ws.Range("N6").Value = 1522
$ws.Range("P6").Value = 112
$ws.Range("Q6").Value = 334
$ws.Range("R6").Value = -164
$ws.Range("S6").Value = -79
$ws.Range("T6").Value = 381
$ws.Range("U6").Value = -47
$ws.Range("V6").Value = 4798
$ws.Range("W6").Value = 3.84
$ws.Range("X6").Value = -1.19
$ws.Range("Y6").Value = -4.89
$ws.Range("Z6").Value = -0.85
$ws.Range("AA6").Value = 334.34
$ws.Range("AB6").Value = 1234.35
$ws.Range("AC6").Value = -297
$ws.Range("AD6").Value = -19.21
$ws.Range("AE6").Value = 5919
$ws.Range("AF6").Value = 0.96
$ws.Range("AG6").Value = 43
$ws.Range("AH6").Value = 0.76
$ws.Range("AI6").Value = -14.47
$ws.Range("AJ6").Value = 25921252

# Clear rows 7, 8, 9 data cells (D:AI), keep A,B,C
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
